$wb = $excel.ActiveWorkbook

# --- subscriptions sheet: insert a new 2020 row at the top of the data ---
$wsSub = $wb.Worksheets.Item("subscriptions")
[void]$wsSub.Rows("2:2").Insert()

$wsSub.Range("A2").Value = 2020
$wsSub.Range("B2").Value = 8
$wsSub.Range("C2").Formula = "=ROUND(B2/`$B`$2*100,1)"
$wsSub.Range("D2").Value = 0
$wsSub.Range("E2").Value = 0

$wsSub.Range("C3").Formula = "=ROUND(B3/`$B`$2*100,1)"
$wsSub.Range("D3").Formula = "=B3-B2"
$wsSub.Range("E3").Formula = "=ROUND(((B3-B2)/B2)*100,1)"

$wsSub.Range("C4").Formula = "=ROUND(B4/`$B`$2*100,1)"
$wsSub.Range("D4").Formula = "=B4-B3"
$wsSub.Range("E4").Formula = "=ROUND(((B4-B3)/B3)*100,1)"

$wsSub.Range("C5").Formula = "=ROUND(B5/`$B`$2*100,1)"
$wsSub.Range("D5").Formula = "=B5-B4"
$wsSub.Range("E5").Formula = "=ROUND(((B5-B4)/B4)*100,1)"

$wsSub.Range("C6").Formula = "=ROUND(B6/`$B`$2*100,1)"
$wsSub.Range("D6").Formula = "=B6-B5"
$wsSub.Range("E6").Formula = "=ROUND(((B6-B5)/B5)*100,1)"

[void]$wsSub.Range("G4").Select()

# --- data source sheet: update notes / drop old url line ---
$wsData = $wb.Worksheets.Item("data source")
$wsData.Range("A1").Value = "Data were collected from Wikipedia and associated links"
[void]$wsData.Range("A2").ClearContents()
[void]$wsData.Range("B5").Select()
